$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate formatting of the last existing data row (row 8) onto the two
# new rows so the new cells pick up the same cell styles (s="1"/"0"/"2")
# used throughout the table.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$ws.Range("A8:F8").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)

# Row 9: new review for com.singleton.strechy / taxi game
$ws.Range("A9").Value = "com.singleton.strechy"
$ws.Range("B9").Value = "taxi game"
$ws.Range("C9").Value = "ctamar115@gmail.com"
$ws.Range("D9").Value = "nirh94846@gmail.com"
$ws.Range("E9").Value = "27/5/2019 15:59"
$ws.Range("F9").Value = "taxi game for every family. Kids and adults. Free car game and works offline too."

# Row 10: another new review for com.singleton.strechy / taxi game
$ws.Range("A10").Value = "com.singleton.strechy"
$ws.Range("B10").Value = "taxi game"
$ws.Range("C10").Value = "syechimovitz@gmail.com"
$ws.Range("D10").Value = "ctamar115@gmail.com"
$ws.Range("E10").Value = "27/5/2019 15:59"
$ws.Range("F10").Value = "I’m a developer and I’m inspired by this game graphics sounds and design. Really really good car game and very creative."

# Update the selected / active cell to match the new last cell of data.
$ws.Range("F10").Select()
